$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Footnote 22 (2nd footnote in the document): drop the stray
#    trailing "a" that was left after "འཁོན། སྣར་ཐང་། པེ་ཅིན།".
# -----------------------------------------------------------------
$fn22 = $d.Footnotes.Item(2)
$t22 = $fn22.Range.Text
if ($t22.EndsWith("a")) {
    $fn22.Range.Text = $t22.Substring(0, $t22.Length - 1)
}

# -----------------------------------------------------------------
# 2. Footnote 31 (the last, 11th, footnote) is empty/bogus - it only
#    contains a single shad "।". Delete the footnote entirely; this
#    also removes its footnoteReference run from the body. Append
#    the folio marker "[༨༤བ]" that used to sit right after it onto
#    the end of the body text, right after the "།།" that precedes it.
# -----------------------------------------------------------------
$fn31 = $d.Footnotes.Item($d.Footnotes.Count)
$fn31.Delete()

# Only touch the very end of the body (a few characters of trailing
# padding is included so the search range does not land exactly on a
# run boundary), so the two other, unrelated, "།།" occurrences that
# appear earlier in the paragraph are left completely untouched.
$tailStart = $d.Content.End - 6
if ($tailStart -lt 0) { $tailStart = 0 }
$tail = $d.Range($tailStart, $d.Content.End)
$tail.Find.Execute("།།", $true, $false, $false, $false, $false, $true, 1, $false, "།།[༨༤བ]", 2)
